# Refresh the crypto price/volume table (GitHub Actions scrape update).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.926.91"
$ws.Range("E2").Value = "  -3.93%  "
$ws.Range("D3").Value = "3.132.53"
$ws.Range("E3").Value = "  -4.39%  "
$ws.Range("E4").Value = "  +0.22%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "606.61"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.79%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.87"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -9.09%  "
$ws.Range("E7").Value = "  +0.31%  "
$ws.Range("D8").Value = "3.118.74"
$ws.Range("E8").Value = "  -4.78%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.518"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -4.73%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.150"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -7.57%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.23"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -8.94%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.469"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -5.80%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000250"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -8.14%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.12"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -9.93%  "
$ws.Range("D15").Value = "3.644.23"
$ws.Range("E15").Value = "  -4.57%  "
$ws.Range("E16").Value = "  +1.48%  "
$ws.Range("D17").Value = "63.940.35"
$ws.Range("E17").Value = "  -3.99%  "
$ws.Range("D18").Value = "3.124.42"
$ws.Range("E18").Value = "  -4.75%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.79"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -8.16%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "477.03"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -5.54%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.60"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -5.81%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.706"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -6.70%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.69"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -5.68%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.51"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -8.53%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "84.07"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.95%  "
$ws.Range("E26").Value = "  -0.15%  "
$ws.Range("E27").Value = "  -8.76%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.41"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -8.62%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.10"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -12.25%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.80"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.78%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.113"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -14.65%  "
$ws.Range("B32").Value = "FirstDigitalUSD"
$ws.Range("C32").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.00"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.03%  "
$ws.Range("B33").Value = "Stacks"
$ws.Range("C33").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.70"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -6.12%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "26.22"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -6.40%  "
$ws.Range("E35").Value = "  -2.62%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.91"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -8.92%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "52.91"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -5.24%  "
$ws.Range("D38").Value = "0.0₃0746"
$ws.Range("E38").Value = "  -6.61%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "457.80"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -9.51%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.92"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -17.33%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0391"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -8.40%  "
$ws.Range("E42").Value = "  -9.67%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.33"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -5.59%  "
$ws.Range("D44").Value = "2.846.54"
$ws.Range("E44").Value = "  -5.63%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.265"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -9.86%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.26"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -13.14%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.42"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.55%  "
$ws.Range("E48").Value = "  +0.00%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "26.03"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -10.29%  "
$ws.Range("E50").Value = "  -5.31%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "118.48"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.45%  "
